$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.738.86"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "'1.635.35"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'212.15"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "'23.48"
$ws.Range("E8").Value = "  +1.99%  "

$ws.Range("D9").Value = "'0.265"
$ws.Range("E9").Value = "  +2.29%  "

$ws.Range("D10").Value = "'0.0613"
$ws.Range("E10").Value = "  +0.19%  "

$ws.Range("D11").Value = "'0.0861"
$ws.Range("E11").Value = "  -3.53%  "

$ws.Range("D12").Value = "'1.866.33"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").Value = "'1.640.93"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").Value = "'0.554"
$ws.Range("E15").Value = "  -1.23%  "

$ws.Range("D16").Value = "'65.18"
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("D17").Value = "'27.686.91"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").Value = "'230.68"
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").Value = "'0.0₃0721"
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").Value = "'10.66"
$ws.Range("E22").Value = "  +4.24%  "

$ws.Range("D23").Value = "'4.37"
$ws.Range("E23").Value = "  +1.52%  "

$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  +3.51%  "

$ws.Range("D25").Value = "'149.21"
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("D26").Value = "'6.90"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("D28").Value = "'15.59"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").Value = "'1.483.78"
$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("E34").Value = "  -0.95%  "

$ws.Range("D35").Value = "'1.55"
$ws.Range("E35").Value = "  -1.21%  "

$ws.Range("E36").Value = "  -1.40%  "

$ws.Range("D37").Value = "'0.960"
$ws.Range("E37").Value = "  +6.60%  "

$ws.Range("D38").Value = "'0.883"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("E39").Value = "  -1.36%  "

$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").Value = "'67.89"
$ws.Range("E43").Value = "  -1.86%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").Value = "'5.34"
$ws.Range("E46").Value = "  -3.96%  "

$ws.Range("D47").Value = "'1.775.57"
$ws.Range("E47").Value = "  -0.27%  "

$ws.Range("D48").Value = "'1.76"
$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").Value = "'87.78"
$ws.Range("E49").Value = "  +1.13%  "

$ws.Range("E50").Value = "  -1.92%  "

$ws.Range("D51").Value = "'0.0991"
$ws.Range("E51").Value = "  -0.22%  "
